$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 37: 504. Base 7
#   columns A-F copy the cell formatting (styles) used by row 3 (A3:F3),
#   column  G copies the "未复习" formatting/value used by row 33 (G33).
# Hyperlinks.Add is issued first so the copied-in formatting (not the
# hyperlink theme) ends up owning the final cell style, matching the
# original author's styling.
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("C37"), "https://leetcode.com/problems/base-7/")

$ws.Range("A3:F3").Copy($ws.Range("A37"))
$ws.Range("G33").Copy($ws.Range("G37"))
$ws.Rows.Item(37).RowHeight = 28

$ws.Range("A37").Value = "504. Base 7"
$ws.Range("B37").Value = "Easy"
$ws.Range("C37").Value = "https://leetcode.com/problems/base-7/"
$ws.Range("D37").Value = 44540
$ws.Range("E37").Value = "数学"
$ws.Range("F37").Value = "简单进制转换"

# ---------------------------------------------------------------------------
# Row 38: 172. Factorial Trailing Zeroes
#   columns A-E copy the formatting used by row 8 (A8:E8),
#   column  F copies the formatting used by row 22 (F22),
#   columns G/H copy the formatting/values used by row 33 (G33:H33).
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("C38"), "https://leetcode.com/problems/factorial-trailing-zeroes/")

$ws.Range("A8:E8").Copy($ws.Range("A38"))
$ws.Range("F22").Copy($ws.Range("F38"))
$ws.Range("G33:H33").Copy($ws.Range("G38"))
$ws.Rows.Item(38).RowHeight = 42

$ws.Range("A38").Value = "172. Factorial Trailing Zeroes"
$ws.Range("B38").Value = "Medium"
$ws.Range("C38").Value = "https://leetcode.com/problems/factorial-trailing-zeroes/"
$ws.Range("D38").Value = 44540
$ws.Range("E38").Value = "数学"
$ws.Range("F38").Value = "n以下所有正整数某个特定质因数的个数求法"
$ws.Range("F38").Characters(2, 19).Font.Name = "宋体"

# ---------------------------------------------------------------------------
# Row 39: 415. Add Strings
#   all columns A-H copy the formatting/values used by row 33 (A33:H33);
#   only the cells that actually differ (A, C, D, F) get overwritten.
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("C39"), "https://leetcode.com/problems/add-strings/")

$ws.Range("A33:H33").Copy($ws.Range("A39"))
$ws.Rows.Item(39).RowHeight = 28

$ws.Range("A39").Value = "415. Add Strings"
$ws.Range("C39").Value = "https://leetcode.com/problems/add-strings/"
$ws.Range("D39").Value = 44542
$ws.Range("F39").Value = "字符串逐位相加；int和char的转换方法"

$ws.Range("F31").Select()

Write-Host "Added rows 37-39"
